# Apply the update described by the commit "chore: update data.xlsx with
# latest information":
#   - Clear out the contents of the "TDB Token" (row 34) and "DRKET Token"
#     (row 35) rows on the only worksheet, leaving the cell formatting
#     intact (this also makes those two strings disappear from the shared
#     strings table when Excel re-saves the workbook).
#   - Update the view state (which cell/range was selected and which cell
#     was scrolled to the top-left) to reflect where the user ended up
#     after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the content (values) of the two rows, keeping number formats/styles.
$ws.Range("A34:K35").ClearContents()

# Update the visible window / selection to match where the user left off.
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("P38").Select()

$wb.Save()
